$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove requisito rows that no longer apply (delete from bottom up so the row
# numbers of not-yet-processed rows stay valid):
#   row 44 -> "LOQ4073 -  Química Geral II  (Requisito)"
#   row 43 -> "LOQ4031 -  Química Geral I  (Requisito)"
#   row 29 -> "LOB1012 -  Estatística  (Requisito)"
$ws.Rows.Item(44).Delete()
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(29).Delete()

# After the three deletions above, the row that used to be row 45
# ("LOQ4095 -  Química Geral Experimental  (Requisito)") is now row 42, the
# last populated row. Insert two new requisito rows right after it.
$ws.Rows.Item(43).Insert()
$ws.Rows.Item(43).RowHeight = 30
$ws.Range("B43").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"
$ws.Range("C43").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n"

$ws.Rows.Item(44).Insert()
$ws.Rows.Item(44).RowHeight = 30
$ws.Range("B44").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)`n"
$ws.Range("C44").Value = "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)`n"

$ws.Range("A1").Select() | Out-Null
